$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 now mirrors the "Ok" cells in B2:E2 - same text and same green
# fill/border formatting (style index 2 in the original workbook).
$ws.Range("F2").Value = $ws.Range("B2").Text
$ws.Range("F2").Interior.Color = $ws.Range("B2").Interior.Color

# Columns C:F are widened to match column B's width.
$ws.Range("C1:F1").ColumnWidth = $ws.Range("B1").ColumnWidth

# The active selection moves from G3 to G8.
$ws.Range("G8").Select()
